# Rename the template placeholders in Intereses.docx to the new (English)
# variable names, matching the supplied OOXML diff:
#   {expediente}  -> {fileNumber}
#   {cantidad}    -> {totalAmount}
#   cantidadInt   -> totalInterests
#   fechaHoy      -> todayDate
#
# For the placeholders whose name is wrapped on its own run (between the
# literal "{" / "}" runs) a temporary bookmark is used to force Word to
# split the run cleanly around the replaced text -- mirroring the
# <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
# wrapped run introduced by the diff (Word inserts those proofErr markers
# itself during its live spell-check pass; that pass is not reachable via
# the COM/object-model surface, so only the run split is reproduced here).

$d = $word.ActiveDocument

function Split-Rename($findText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $bmName = "tmpRenameBm"
    $d.Bookmarks.Add($bmName, $rng)
    $bmRange = $d.Bookmarks($bmName).Range
    $bmRange.Text = $newText
    $d.Bookmarks($bmName).Delete()
}

# 1) "ID {expediente} " -> "ID {" / "fileNumber" / "} "
Split-Rename "expediente" "fileNumber"

# 2) "{cantidad} " -> "{" / "totalAmount" / "} "
$rng = $d.Content
$rng.Find.Execute("{cantidad}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordRng = $d.Range($rng.Start + 1, $rng.End - 1)
$d.Bookmarks.Add("tmpRenameBm2", $wordRng)
$bm2 = $d.Bookmarks("tmpRenameBm2").Range
$bm2.Text = "totalAmount"
$d.Bookmarks("tmpRenameBm2").Delete()

# 3) "cantidadInt" -> "totalInterests" (already an isolated run; plain replace)
$rng = $d.Content
$rng.Find.Execute("cantidadInt", $true, $false, $false, $false, $false, $true, 1, $false, "totalInterests", 2)

# 4) "fechaHoy" -> "todayDate" (already an isolated run; plain replace)
$rng = $d.Content
$rng.Find.Execute("fechaHoy", $true, $false, $false, $false, $false, $true, 1, $false, "todayDate", 2)

Write-Output "Renames applied"
